$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8702740669250488
$ws.Range("B1").Value = 3.170135259628296
$ws.Range("C1").Value = 2.943016290664673
$ws.Range("D1").Value = 1.671888828277588
$ws.Range("E1").Value = 1.284761428833008
